$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B11 currently holds the text "R40"; change it to the text "1".
# Leading apostrophe forces Excel to store it as text (not a number),
# matching the workbook's existing shared-string ("t=s") cell type.
$ws.Range("B11").Value = "'1"
